# Update Notion sync timestamps and the two number columns that differ
# between the Mac/Win exports ("Đầy đủ" / AC and "Tổng công" / AF).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9 were last edited at 12:51 on 2024-07-19 -> now 13:34 on 2024-07-20
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-20T13:34:00.000Z"
}

# Rows 10-20 were last edited at 12:52 on 2024-07-19 -> now 13:35 on 2024-07-20
for ($r = 10; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-07-20T13:35:00.000Z"
}

# properties.Đầy đủ.number (col AC = 29) and properties.Tổng công.number (col AF = 32)
# bump by 1 on the rows whose counts changed.
$ws.Cells.Item(8, 29).Value = 18
$ws.Cells.Item(8, 32).Value = 19

$ws.Cells.Item(14, 29).Value = 18
$ws.Cells.Item(14, 32).Value = 18.5

$ws.Cells.Item(17, 29).Value = 18
$ws.Cells.Item(17, 32).Value = 18.5

$ws.Cells.Item(19, 29).Value = 18
$ws.Cells.Item(19, 32).Value = 19

$ws.Cells.Item(20, 29).Value = 19
$ws.Cells.Item(20, 32).Value = 19
